$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.008725
$ws.Range("E2").Value = -0.0459
$ws.Range("F2").Value = 0.05435
$ws.Range("I2").Value = 0.00001249604769928682
$ws.Range("J2").Value = 0.0000107858846640895
$ws.Range("K2").Value = 7809
$ws.Range("L2").Value = 0.5292801952013013
$ws.Range("M2").Value = 3691.3
$ws.Range("N2").Value = 0.04937123994532289
$ws.Range("O2").Value = 0.4726981687796133
$ws.Range("P2").Value = 3387.7
$ws.Range("Q2").Value = 0.04531058151945665
$ws.Range("R2").Value = 0.4338199513381995
$ws.Range("S2").Value = 303.6
$ws.Range("T2").Value = 0.08224744669899493
$ws.Range("U2").Value = 58686.4
$ws.Range("V2").Value = 0.7849322287343746
$ws.Range("W2").Value = 0.070630868034638
$ws.Range("X2").Value = 0.04817607143282353
$ws.Range("Y2").Value = 0.02245479660181447
$ws.Range("Z2").Value = 0.2926216378263702
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03784420390835442
$ws.Range("AC2").Value = -0.03784420390835442
$ws.Range("AD2").Value = 19969
$ws.Range("AE2").Value = 0.2281665612236116
$ws.Range("AF2").Value = 19969.22816656122
$ws.Range("AG2").Value = -38717.17183343878
$ws.Range("AH2").Value = 0.2107894433268605
$ws.Range("AI2").Value = 0.1849445650142841
$ws.Range("AJ2").Value = -1.074014302259402
$ws.Range("AK2").Value = -0.7855337073479899
$ws.Range("AN2").Value = 86821.73913043478
$ws.Range("AP2").Value = -168335.5297106034

# Row 3
$ws.Range("B3").Value = "Hang Seng Bank Limited (SEHK:11)"
$ws.Range("D3").Value = -0.0241
$ws.Range("E3").Value = -0.0532
$ws.Range("F3").ClearContents()
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2622.7
$ws.Range("L3").Value = 0.6107113750145534
$ws.Range("M3").Value = 1999
$ws.Range("N3").Value = 0.06063418688311768
$ws.Range("O3").Value = 0.7621916345750563
$ws.Range("P3").Value = 1999
$ws.Range("Q3").Value = 0.06063418688311768
$ws.Range("R3").Value = 0.7621916345750563
$ws.Range("U3").Value = 12631.6
$ws.Range("V3").Value = 0.3831449700013953
$ws.Range("W3").Value = 0.1186071226681741
$ws.Range("X3").Value = 0.03808546457284634
$ws.Range("Y3").Value = 0.08052165809532776
$ws.Range("Z3").Value = 0.231369739025494
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03563832963014231
$ws.Range("AC3").Value = -0.03563832963014231
$ws.Range("AD3").Value = 5211.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 5211.7
$ws.Range("AG3").Value = -7419.900000000001
$ws.Range("AH3").Value = 0.1365037624509231
$ws.Range("AI3").Value = 0.1865659085946254
$ws.Range("AJ3").Value = -0.2904263688777727
$ws.Range("AK3").Value = -0.4848562074846602
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# Row 4
$ws.Range("B4").Value = "BOC Hong Kong (Holdings) Limited (SEHK:2388)"
$ws.Range("D4").Value = 0.0852
$ws.Range("E4").Value = 0.04019999999999999
$ws.Range("F4").Value = -0.0293
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4067.2
$ws.Range("L4").Value = 0.5688391608391609
$ws.Range("M4").Value = 922
$ws.Range("N4").Value = 0.02877141082890996
$ws.Range("O4").Value = 0.2266915814319434
$ws.Range("P4").Value = 922
$ws.Range("Q4").Value = 0.02877141082890996
$ws.Range("R4").Value = 0.2266915814319434
$ws.Range("U4").Value = 26819.3
$ws.Range("V4").Value = 0.83690791588263
$ws.Range("W4").Value = 0.108868593209632
$ws.Range("X4").Value = 0.03941424029999772
$ws.Range("Y4").Value = 0.06945435290963428
$ws.Range("Z4").Value = 0.3936574354456862
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0359936505297057
$ws.Range("AC4").Value = -0.0359936505297057
$ws.Range("AD4").Value = 7071.3
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 7071.3
$ws.Range("AG4").Value = -19748
$ws.Range("AH4").Value = 0.1807730654191272
$ws.Range("AI4").Value = 0.1490059317480219
$ws.Range("AJ4").Value = -1.605828732201956
$ws.Range("AK4").Value = -0.9569127594828757
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("B5").Value = "Dah Sing Banking Group Limited (SEHK:2356)"
$ws.Range("D5").Value = 0.0142
$ws.Range("E5").Value = -0.0386
$ws.Range("F5").ClearContents()
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 236.3
$ws.Range("L5").Value = 0.4256123919308357
$ws.Range("M5").Value = 92.40000000000001
$ws.Range("N5").Value = 0.06402439024390244
$ws.Range("O5").Value = 0.3910283537875582
$ws.Range("P5").Value = 92.40000000000001
$ws.Range("Q5").Value = 0.06402439024390244
$ws.Range("R5").Value = 0.3910283537875582
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 3297.3
$ws.Range("V5").Value = 2.284714523281596
$ws.Range("W5").Value = 0.06690828779341394
$ws.Range("X5").Value = 0.04814081399830369
$ws.Range("Y5").Value = 0.01876747379511025
$ws.Range("Z5").Value = 0.2382423618262959
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03743737888972589
$ws.Range("AC5").Value = -0.03743737888972589
$ws.Range("AD5").Value = 911.6
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 911.6
$ws.Range("AG5").Value = -2385.7
$ws.Range("AH5").Value = 0.3871241719041957
$ws.Range("AI5").Value = 0.2012806359019651
$ws.Range("AJ5").Value = 2.531246684350132
$ws.Range("AK5").Value = -1.936916456929448
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# Row 6
$ws.Range("B6").Value = "Dah Sing Financial Holdings Limited (SEHK:440)"
$ws.Range("D6").Value = 0.00325
$ws.Range("E6").Value = -0.0562
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 172.7
$ws.Range("L6").Value = 0.299047619047619
$ws.Range("M6").Value = 60.6
$ws.Range("N6").Value = 0.06728847435043304
$ws.Range("O6").Value = 0.3508975101331789
$ws.Range("P6").Value = 60.6
$ws.Range("Q6").Value = 0.06728847435043304
$ws.Range("R6").Value = 0.3508975101331789
$ws.Range("U6").Value = 1612.2
$ws.Range("V6").Value = 1.790139906728847
$ws.Range("W6").Value = 0.05158303464755077
$ws.Range("X6").Value = 0.05625434143074226
$ws.Range("Y6").Value = -0.004671306783191491
$ws.Range("Z6").Value = 0.1717114652711703
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.03825102892698295
$ws.Range("AC6").Value = -0.03825102892698295
$ws.Range("AD6").Value = 913
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 913
$ws.Range("AG6").Value = -699.2
$ws.Range("AH6").Value = 0.5034186149095722
$ws.Range("AI6").Value = 0.1734257764270111
$ws.Range("AJ6").Value = -3.471698113207548
$ws.Range("AK6").Value = -0.1914410097746626
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()

# Row 7
$ws.Range("B7").Value = "The Bank of East Asia, Limited (SEHK:23)"
$ws.Range("D7").Value = -0.0377
$ws.Range("E7").Value = -0.1
$ws.Range("F7").Value = 0.138
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 489.3
$ws.Range("L7").Value = 0.2882644043831742
$ws.Range("M7").Value = 219.6
$ws.Range("N7").Value = 0.0352521912222686
$ws.Range("O7").Value = 0.4488044144696505
$ws.Range("P7").Value = 219.6
$ws.Range("Q7").Value = 0.0352521912222686
$ws.Range("R7").Value = 0.4488044144696505
$ws.Range("U7").Value = 11898.3
$ws.Range("V7").Value = 1.910023437249173
$ws.Range("W7").Value = 0.03741540814375836
$ws.Range("X7").Value = 0.04821132886734338
$ws.Range("Y7").Value = -0.01079592072358501
$ws.Range("Z7").Value = 0.4113114277406223
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.0385594978466537
$ws.Range("AC7").Value = -0.0385594978466537
$ws.Range("AD7").Value = 3955.5
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 3955.5
$ws.Range("AG7").Value = -7942.799999999999
$ws.Range("AH7").Value = 0.3883690561517541
$ws.Range("AI7").Value = 0.2227083087005726
$ws.Range("AJ7").Value = 4.635695109139723
$ws.Range("AK7").Value = -1.354825504042575
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# Row 8
$ws.Range("B8").Value = "Chong Hing Bank Limited (SEHK:1111)"
$ws.Range("D8").Value = 0.104
$ws.Range("E8").Value = 0.0863
$ws.Range("I8").Value = 0.0003845779886426318
$ws.Range("J8").Value = 0.0003211604383218347
$ws.Range("K8").Value = 220.8
$ws.Range("L8").Value = 0.460575719649562
$ws.Range("M8").Value = 397.7
$ws.Range("N8").Value = 0.3372911542702062
$ws.Range("O8").Value = 1.801177536231884
$ws.Range("P8").Value = 94.09999999999999
$ws.Range("Q8").Value = 0.07980663217708422
$ws.Range("R8").Value = 0.426177536231884
$ws.Range("S8").Value = 303.6
$ws.Range("T8").Value = 0.7633894895649987
$ws.Range("U8").Value = 2427.7
$ws.Range("V8").Value = 2.058943261809855
$ws.Range("W8").Value = 0.07435344827586207
$ws.Range("X8").Value = 0.06905424067567952
$ws.Range("Y8").Value = 0.005299207600182557
$ws.Range("Z8").Value = 0.1237015029955495
$ws.Range("AA8").Value = 0.00003972802892312043
$ws.Range("AB8").Value = 0.0390514625983992
$ws.Range("AC8").Value = -0.03901173456947608
$ws.Range("AD8").Value = 1905.9
$ws.Range("AE8").Value = 0.2281665612236116
$ws.Range("AF8").Value = 1906.128166561224
$ws.Range("AG8").Value = -521.5718334387761
$ws.Range("AH8").Value = 0.6178240517899143
$ws.Range("AI8").Value = 0.3790779168386673
$ws.Range("AJ8").Value = -0.7932311647826749
$ws.Range("AK8").Value = -0.2005560964635878
$ws.Range("AN8").Value = 8286.521739130434
$ws.Range("AP8").Value = -2267.703623646853
